$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(326).Insert()

$ws.Cells.Item(326, 1).Value = 9
$ws.Cells.Item(326, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(326, 3).Value = "Metropolitana"
$ws.Cells.Item(326, 4).Value = 45131
$ws.Cells.Item(326, 5).Value = 13
$ws.Cells.Item(326, 6).Value = 100112026
$ws.Cells.Item(326, 7).Value = "Haba"
$ws.Cells.Item(326, 8).Value = "Sin especificar"
$ws.Cells.Item(326, 9).Value = "Primera"
$ws.Cells.Item(326, 10).Value = 52
$ws.Cells.Item(326, 11).Value = 15000
$ws.Cells.Item(326, 12).Value = 16000
$ws.Cells.Item(326, 13).Value = 15500
$ws.Cells.Item(326, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(326, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(326, 16).Value = 620
$ws.Cells.Item(326, 17).Value = 25
$ws.Cells.Item(326, 18).Value = "Hortaliza"
